# Commit for Dietitian POST phase1
#
# Adds a new "DietitianPost" worksheet (with a header row + one data row of
# dietitian-registration test fixtures, an Email hyperlink and matching
# column widths) after the existing "AdminLogin" sheet, and updates the
# selection left behind on "AdminLogin".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet            # "AdminLogin"

# --- 1. AdminLogin: move the saved selection from D12 to K6 -----------------
[void]$ws1.Range("K6").Select()

# --- 2. Insert the new sheet right after AdminLogin --------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "DietitianPost"

# --- 3. Header row (row 1) ---------------------------------------------------
$ws2.Cells.Item(1,1).Value  = "ContactNumber"
$ws2.Cells.Item(1,3).Value  = "Education"
$ws2.Cells.Item(1,4).Value  = "Email"
$ws2.Cells.Item(1,5).Value  = "Firstname"
$ws2.Cells.Item(1,6).Value  = "HospitalCity"
$ws2.Cells.Item(1,7).Value  = "HospitalName"
$ws2.Cells.Item(1,8).Value  = "HospitalPincode"
$ws2.Cells.Item(1,9).Value  = "HospitalStreet"
$ws2.Cells.Item(1,10).Value = "Lastname"

# --- 4. Data row (row 2) - static text fields --------------------------------
$ws2.Cells.Item(2,3).Value  = "Mphars"
$ws2.Cells.Item(2,6).Value  = "Hartford"
$ws2.Cells.Item(2,7).Value  = "Saintfrancis"
$ws2.Cells.Item(2,9).Value  = "Pinwheelstreet"
$ws2.Cells.Item(2,10).Value = "Antonyrt"

# --- 5. Data row (row 2) - numeric fields ------------------------------------
$ws2.Cells.Item(2,1).Value = 8807306309
$ws2.Cells.Item(2,8).Value = 160741

# --- 6. DateOfBirth column + remaining dynamic fields ------------------------
$ws2.Cells.Item(1,2).Value = "DateOfBirth"
$ws2.Cells.Item(2,2).Value = "2024-07-26T18:14:08.570Z"
$ws2.Cells.Item(2,5).Value = "Markty"
$ws2.Cells.Item(2,4).Value = "shr10@gmail.com"

# --- 7. Email hyperlink (D2) with the Hyperlink cell style -------------------
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:shr10@gmail.com")
$ws2.Cells.Item(2,4).Style = "Hyperlink"

# --- 8. Column widths (best-fit on the populated columns) -------------------
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(3).AutoFit()
$ws2.Columns.Item(4).AutoFit()
$ws2.Columns.Item(6).AutoFit()
$ws2.Columns.Item(7).AutoFit()
$ws2.Columns.Item(8).AutoFit()
$ws2.Columns.Item(9).AutoFit()
$ws2.Columns.Item(2).ColumnWidth = $ws2.Columns.Item(1).ColumnWidth

# --- 9. Leave the new sheet selected/active, as in the saved file -----------
[void]$ws2.Range("D2").Select()
$ws2.Activate()
